# Update "想去人数" (want-to-go count) figures on both the "展览" sheet
# and the "全部类型" sheet, each incrementing by 1:
#   F9:  4497 -> 4498
#   F10: 4373 -> 4374

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F9").Value = 4498
    $ws.Range("F10").Value = 4374
}
